$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 814  # H98: was 598
$ws.Cells.Item(98, 9).Value = 814  # I98: was 559
$ws.Cells.Item(98, 10).Value = 0  # J98: was 695.5
$ws.Cells.Item(98, 11).Value = 814  # K98: was 559
$ws.Cells.Item(98, 12).Value = 0  # L98: was 695.5
$ws.Cells.Item(98, 13).Value = 684  # M98: was 939
$ws.Cells.Item(98, 14).ClearContents()  # N98: was -3691.5

$ws.Cells.Item(112, 8).Value = 1289.3823  # H112: was 1378.8276
$ws.Cells.Item(112, 10).Value = 1313.303  # J112: was 1410.2142
$ws.Cells.Item(112, 12).Value = 3939.909000000001  # L112: was 4230.642599999999
$ws.Cells.Item(112, 14).Value = -6155.909000000001  # N112: was -6446.642599999999

$ws.Cells.Item(122, 8).Value = 814  # H122: was 598
$ws.Cells.Item(122, 9).Value = 814  # I122: was 559
$ws.Cells.Item(122, 10).Value = 0  # J122: was 695.5
$ws.Cells.Item(122, 11).Value = 2442  # K122: was 1677
$ws.Cells.Item(122, 12).Value = 0  # L122: was 2086.5
$ws.Cells.Item(122, 13).Value = 8  # M122: was 773
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -6986.5

$ws.Cells.Item(123, 8).Value = 32970.855  # H123: was 36799.668
$ws.Cells.Item(123, 10).Value = 32970.855  # J123: was 36799.668
$ws.Cells.Item(123, 12).Value = 32970.855  # L123: was 36799.668
$ws.Cells.Item(123, 14).Value = -42770.855  # N123: was -46599.668

$ws.Cells.Item(136, 8).Value = 58401.43  # H136: was 58975
$ws.Cells.Item(136, 10).Value = 58401.43  # J136: was 58975
$ws.Cells.Item(136, 12).Value = 58401.43  # L136: was 58975
$ws.Cells.Item(136, 14).Value = -68601.42999999999  # N136: was -69175

$ws.Cells.Item(138, 8).Value = 3313.9019  # H138: was 3442.06
$ws.Cells.Item(138, 9).Value = 3375  # I138: was 3714.2856
$ws.Cells.Item(138, 10).Value = 3302.535  # J138: was 3397.7441
$ws.Cells.Item(138, 11).Value = 10125  # K138: was 11142.8568
$ws.Cells.Item(138, 12).Value = 9907.605  # L138: was 10193.2323
$ws.Cells.Item(138, 13).Value = -4985  # M138: was -6002.856800000001
$ws.Cells.Item(138, 14).Value = -20187.605  # N138: was -20473.2323

$ws.Cells.Item(139, 8).Value = 69875.60000000001  # H139: was 69924.5
$ws.Cells.Item(139, 10).Value = 69875.60000000001  # J139: was 69924.5
$ws.Cells.Item(139, 12).Value = 69875.60000000001  # L139: was 69924.5
$ws.Cells.Item(139, 14).Value = -80155.60000000001  # N139: was -80204.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 24505.023  # H32: was 26105.293
$ws.Cells.Item(32, 9).Value = 6889.8076  # I32: was 7274.0547
$ws.Cells.Item(32, 10).Value = 149412.9  # J32: was 178847.56
$ws.Cells.Item(32, 11).Value = 6889.8076  # K32: was 7274.0547
$ws.Cells.Item(32, 12).Value = 149412.9  # L32: was 178847.56
$ws.Cells.Item(32, 13).Value = -6602.8076  # M32: was -6987.0547
$ws.Cells.Item(32, 14).Value = -149986.9  # N32: was -179421.56

$ws.Cells.Item(61, 8).Value = 1715.7333  # H61: was 1636.8077
$ws.Cells.Item(61, 9).Value = 1139.3182  # I61: was 1112.1052
$ws.Cells.Item(61, 10).Value = 3300.875  # J61: was 3061
$ws.Cells.Item(61, 11).Value = 1139.3182  # K61: was 1112.1052
$ws.Cells.Item(61, 12).Value = 3300.875  # L61: was 3061
$ws.Cells.Item(61, 13).Value = -927.3181999999999  # M61: was -900.1052
$ws.Cells.Item(61, 14).Value = -3724.875  # N61: was -3485

$ws.Cells.Item(136, 8).Value = 1715.7333  # H136: was 1636.8077
$ws.Cells.Item(136, 9).Value = 1139.3182  # I136: was 1112.1052
$ws.Cells.Item(136, 10).Value = 3300.875  # J136: was 3061
$ws.Cells.Item(136, 11).Value = 3417.9546  # K136: was 3336.3156
$ws.Cells.Item(136, 12).Value = 9902.625  # L136: was 9183
$ws.Cells.Item(136, 13).Value = -867.9546  # M136: was -786.3155999999999
$ws.Cells.Item(136, 14).Value = -15002.625  # N136: was -14283

$ws.Cells.Item(139, 8).Value = 39930  # H139: was 53395
$ws.Cells.Item(139, 10).Value = 39930  # J139: was 53395
$ws.Cells.Item(139, 12).Value = 39930  # L139: was 53395
$ws.Cells.Item(139, 14).Value = -50210  # N139: was -63675

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 78721.11  # H138: was 81348.75
$ws.Cells.Item(138, 10).Value = 78721.11  # J138: was 81348.75
$ws.Cells.Item(138, 12).Value = 78721.11  # L138: was 81348.75
$ws.Cells.Item(138, 14).Value = -89001.11  # N138: was -91628.75

$ws.Cells.Item(140, 8).Value = 49050  # H140: was 49562.5
$ws.Cells.Item(140, 10).Value = 49050  # J140: was 49562.5
$ws.Cells.Item(140, 12).Value = 49050  # L140: was 49562.5
$ws.Cells.Item(140, 14).Value = -59410  # N140: was -59922.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 983.73914  # H107: was 1043
$ws.Cells.Item(107, 9).Value = 1082.6875  # I107: was 1148.8667
$ws.Cells.Item(107, 10).Value = 757.5714  # J107: was 816.1429000000001
$ws.Cells.Item(107, 11).Value = 1082.6875  # K107: was 1148.8667
$ws.Cells.Item(107, 12).Value = 757.5714  # L107: was 816.1429000000001
$ws.Cells.Item(107, 13).Value = 837.3125  # M107: was 771.1333
$ws.Cells.Item(107, 14).Value = -4597.5714  # N107: was -4656.1429

$ws.Cells.Item(132, 8).Value = 65220644  # H132: was 65220584
$ws.Cells.Item(132, 9).Value = 83337470  # I132: was 83337390
$ws.Cells.Item(132, 10).Value = 45456830  # J132: was 45456790
$ws.Cells.Item(132, 11).Value = 250012410  # K132: was 250012170
$ws.Cells.Item(132, 12).Value = 136370490  # L132: was 136370370
$ws.Cells.Item(132, 13).Value = -250009880  # M132: was -250009640
$ws.Cells.Item(132, 14).Value = -136375550  # N132: was -136375430

$ws.Cells.Item(134, 8).Value = 1485.25  # H134: was 1523.625
$ws.Cells.Item(134, 9).Value = 886.36365  # I134: was 898.1667
$ws.Cells.Item(134, 10).Value = 2802.8  # J134: was 3400
$ws.Cells.Item(134, 11).Value = 2659.09095  # K134: was 2694.5001
$ws.Cells.Item(134, 12).Value = 8408.400000000001  # L134: was 10200
$ws.Cells.Item(134, 13).Value = -124.0909499999998  # M134: was -159.5001000000002
$ws.Cells.Item(134, 14).Value = -13478.4  # N134: was -15270

$ws.Cells.Item(138, 8).Value = 74561.25  # H138: was 72479
$ws.Cells.Item(138, 10).Value = 74561.25  # J138: was 72479
$ws.Cells.Item(138, 12).Value = 74561.25  # L138: was 72479
$ws.Cells.Item(138, 14).Value = -84841.25  # N138: was -82759

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 13.833333  # H2: was 14.1
$ws.Cells.Item(2, 10).Value = 13.888889  # J2: was 14.285714
$ws.Cells.Item(2, 12).Value = 83.33333400000001  # L2: was 85.71428400000001
$ws.Cells.Item(2, 14).Value = -309.333334  # N2: was -311.714284

$ws.Cells.Item(12, 8).Value = 41.22222  # H12: was 43.588234
$ws.Cells.Item(12, 10).Value = 46.333332  # J12: was 49.57143
$ws.Cells.Item(12, 12).Value = 138.999996  # L12: was 148.71429
$ws.Cells.Item(12, 14).Value = -484.999996  # N12: was -494.71429

$ws.Cells.Item(21, 8).Value = 500  # H21: was 600
$ws.Cells.Item(21, 10).Value = 500  # J21: was 600
$ws.Cells.Item(21, 12).Value = 1500  # L21: was 1800
$ws.Cells.Item(21, 14).Value = -1846  # N21: was -2146

$ws.Cells.Item(26, 8).Value = 1164.9333  # H26: was 1094.3125
$ws.Cells.Item(26, 9).Value = 241.42857  # I26: was 215.625
$ws.Cells.Item(26, 11).Value = 724.28571  # K26: was 646.875
$ws.Cells.Item(26, 13).Value = -436.28571  # M26: was -358.875

$ws.Cells.Item(39, 8).Value = 1900  # H39: was 1875
$ws.Cells.Item(39, 10).Value = 2250  # J39: was 2333.3333
$ws.Cells.Item(39, 12).Value = 6750  # L39: was 6999.999899999999
$ws.Cells.Item(39, 14).Value = -7338  # N39: was -7587.999899999999

$ws.Cells.Item(57, 8).Value = 2250  # H57: was 4225
$ws.Cells.Item(57, 9).Value = 500  # I57: was 0
$ws.Cells.Item(57, 10).Value = 4000  # J57: was 4225
$ws.Cells.Item(57, 11).Value = 1500  # K57: was 0
$ws.Cells.Item(57, 12).Value = 12000  # L57: was 12675
$ws.Cells.Item(57, 13).Value = -941  # M57: was None
$ws.Cells.Item(57, 14).Value = -13118  # N57: was -13793

$ws.Cells.Item(58, 8).Value = 1600  # H58: was 1633.3334
$ws.Cells.Item(58, 10).Value = 1600  # J58: was 1633.3334
$ws.Cells.Item(58, 12).Value = 4800  # L58: was 4900.0002
$ws.Cells.Item(58, 14).Value = -5056  # N58: was -5156.0002

$ws.Cells.Item(62, 8).Value = 3800  # H62: was 3480.8462
$ws.Cells.Item(62, 9).Value = 3800  # I62: was 2453
$ws.Cells.Item(62, 10).Value = 0  # J62: was 3937.6667
$ws.Cells.Item(62, 11).Value = 11400  # K62: was 7359
$ws.Cells.Item(62, 12).Value = 0  # L62: was 11813.0001
$ws.Cells.Item(62, 13).Value = -10714  # M62: was -6673
$ws.Cells.Item(62, 14).ClearContents()  # N62: was -13185.0001

$ws.Cells.Item(65, 8).Value = 3800  # H65: was 3480.8462
$ws.Cells.Item(65, 9).Value = 3800  # I65: was 2453
$ws.Cells.Item(65, 10).Value = 0  # J65: was 3937.6667
$ws.Cells.Item(65, 11).Value = 34200  # K65: was 22077
$ws.Cells.Item(65, 12).Value = 0  # L65: was 35439.0003
$ws.Cells.Item(65, 13).Value = -30768  # M65: was -18645
$ws.Cells.Item(65, 14).ClearContents()  # N65: was -42303.0003

$ws.Cells.Item(75, 8).Value = 2026.5  # H75: was 2457.5454
$ws.Cells.Item(75, 9).Value = 862.8  # I75: was 971
$ws.Cells.Item(75, 10).Value = 2555.4546  # J75: was 3015
$ws.Cells.Item(75, 11).Value = 2588.4  # K75: was 2913
$ws.Cells.Item(75, 12).Value = 7666.3638  # L75: was 9045
$ws.Cells.Item(75, 13).Value = -1590.4  # M75: was -1915
$ws.Cells.Item(75, 14).Value = -9662.363799999999  # N75: was -11041

$ws.Cells.Item(78, 8).Value = 2026.5  # H78: was 2457.5454
$ws.Cells.Item(78, 9).Value = 862.8  # I78: was 971
$ws.Cells.Item(78, 10).Value = 2555.4546  # J78: was 3015
$ws.Cells.Item(78, 11).Value = 7765.2  # K78: was 8739
$ws.Cells.Item(78, 12).Value = 22999.0914  # L78: was 27135
$ws.Cells.Item(78, 13).Value = -2773.2  # M78: was -3747
$ws.Cells.Item(78, 14).Value = -32983.0914  # N78: was -37119

$ws.Cells.Item(88, 8).Value = 0  # H88: was 3000
$ws.Cells.Item(88, 10).Value = 0  # J88: was 3000
$ws.Cells.Item(88, 12).Value = 0  # L88: was 9000
$ws.Cells.Item(88, 14).ClearContents()  # N88: was -9856

$ws.Cells.Item(91, 8).Value = 0  # H91: was 3000
$ws.Cells.Item(91, 10).Value = 0  # J91: was 3000
$ws.Cells.Item(91, 12).Value = 0  # L91: was 9000
$ws.Cells.Item(91, 14).ClearContents()  # N91: was -11964

$ws.Cells.Item(103, 8).Value = 668.1818  # H103: was 1366.0714
$ws.Cells.Item(103, 9).Value = 668.1818  # I103: was 1082.5
$ws.Cells.Item(103, 10).Value = 0  # J103: was 1578.75
$ws.Cells.Item(103, 11).Value = 2004.5454  # K103: was 3247.5
$ws.Cells.Item(103, 12).Value = 0  # L103: was 4736.25
$ws.Cells.Item(103, 13).Value = -1125.5454  # M103: was -2368.5
$ws.Cells.Item(103, 14).ClearContents()  # N103: was -6494.25

$ws.Cells.Item(122, 8).Value = 608.25  # H122: was 762.25
$ws.Cells.Item(122, 9).Value = 359.2  # I122: was 797
$ws.Cells.Item(122, 10).Value = 786.1429000000001  # J122: was 757.2857
$ws.Cells.Item(122, 11).Value = 3232.8  # K122: was 7173
$ws.Cells.Item(122, 12).Value = 7075.2861  # L122: was 6815.571300000001
$ws.Cells.Item(122, 13).Value = -782.7999999999997  # M122: was -4723
$ws.Cells.Item(122, 14).Value = -11975.2861  # N122: was -11715.5713

$ws.Cells.Item(125, 8).Value = 1690  # H125: was 1247.5
$ws.Cells.Item(125, 9).Value = 1225  # I125: was 1247.5
$ws.Cells.Item(125, 10).Value = 2000  # J125: was 0
$ws.Cells.Item(125, 11).Value = 3675  # K125: was 3742.5
$ws.Cells.Item(125, 12).Value = 6000  # L125: was 0
$ws.Cells.Item(125, 13).Value = 1245  # M125: was 1177.5
$ws.Cells.Item(125, 14).Value = -15840  # N125: was None

$ws.Cells.Item(131, 8).Value = 13481.862  # H131: was 17515.895
$ws.Cells.Item(131, 9).Value = 453.25  # I131: was 479.9091
$ws.Cells.Item(131, 10).Value = 15566.44  # J131: was 20923.092
$ws.Cells.Item(131, 11).Value = 1359.75  # K131: was 1439.7273
$ws.Cells.Item(131, 12).Value = 46699.32  # L131: was 62769.276
$ws.Cells.Item(131, 13).Value = 3680.25  # M131: was 3600.2727
$ws.Cells.Item(131, 14).Value = -56779.32  # N131: was -72849.276

$ws.Cells.Item(132, 8).Value = 478548.94  # H132: was 402033.9
$ws.Cells.Item(132, 9).Value = 681  # I132: was 553.36365
$ws.Cells.Item(132, 11).Value = 6129  # K132: was 4980.27285
$ws.Cells.Item(132, 13).Value = -3599  # M132: was -2450.27285

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 83421800  # H80: was 143006800
$ws.Cells.Item(80, 9).Value = 143006140  # I80: was 166840160
$ws.Cells.Item(80, 10).Value = 3733.2  # J80: was 6666
$ws.Cells.Item(80, 11).Value = 143006140  # K80: was 166840160
$ws.Cells.Item(80, 12).Value = 3733.2  # L80: was 6666
$ws.Cells.Item(80, 13).Value = -143005142  # M80: was -166839162
$ws.Cells.Item(80, 14).Value = -5729.2  # N80: was -8662

$ws.Cells.Item(83, 8).Value = 83421800  # H83: was 143006800
$ws.Cells.Item(83, 9).Value = 143006140  # I83: was 166840160
$ws.Cells.Item(83, 10).Value = 3733.2  # J83: was 6666
$ws.Cells.Item(83, 11).Value = 715030700  # K83: was 834200800
$ws.Cells.Item(83, 12).Value = 18666  # L83: was 33330
$ws.Cells.Item(83, 13).Value = -715025708  # M83: was -834195808
$ws.Cells.Item(83, 14).Value = -28650  # N83: was -43314

$ws.Cells.Item(122, 8).Value = 1614.4242  # H122: was 1725.2122
$ws.Cells.Item(122, 9).Value = 1325.7693  # I122: was 1440.2609
$ws.Cells.Item(122, 10).Value = 2686.5715  # J122: was 2380.6
$ws.Cells.Item(122, 11).Value = 3977.3079  # K122: was 4320.7827
$ws.Cells.Item(122, 12).Value = 8059.7145  # L122: was 7141.799999999999
$ws.Cells.Item(122, 13).Value = -1527.3079  # M122: was -1870.7827
$ws.Cells.Item(122, 14).Value = -12959.7145  # N122: was -12041.8

$ws.Cells.Item(140, 8).Value = 115376.664  # H140: was 134120
$ws.Cells.Item(140, 10).Value = 115376.664  # J140: was 134120
$ws.Cells.Item(140, 12).Value = 115376.664  # L140: was 134120
$ws.Cells.Item(140, 14).Value = -125736.664  # N140: was -144480

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2429.8096  # H7: was 2643.3333
$ws.Cells.Item(7, 9).Value = 1622.1818  # I7: was 1686.4
$ws.Cells.Item(7, 10).Value = 3318.2  # J7: was 3839.5
$ws.Cells.Item(7, 11).Value = 1622.1818  # K7: was 1686.4
$ws.Cells.Item(7, 12).Value = 3318.2  # L7: was 3839.5
$ws.Cells.Item(7, 13).Value = -1510.1818  # M7: was -1574.4
$ws.Cells.Item(7, 14).Value = -3542.2  # N7: was -4063.5

$ws.Cells.Item(40, 8).Value = 44279.707  # H40: was 42504.68
$ws.Cells.Item(40, 9).Value = 85612  # I40: was 68842.664
$ws.Cells.Item(40, 10).Value = 2947.4167  # J40: was 2997.7
$ws.Cells.Item(40, 11).Value = 85612  # K40: was 68842.664
$ws.Cells.Item(40, 12).Value = 2947.4167  # L40: was 2997.7
$ws.Cells.Item(40, 13).Value = -85476  # M40: was -68706.664
$ws.Cells.Item(40, 14).Value = -3219.4167  # N40: was -3269.7

$ws.Cells.Item(122, 8).Value = 3219.476  # H122: was 3159.682
$ws.Cells.Item(122, 9).Value = 3034.1333  # I122: was 2957.25
$ws.Cells.Item(122, 10).Value = 3682.8333  # J122: was 3699.5
$ws.Cells.Item(122, 11).Value = 9102.3999  # K122: was 8871.75
$ws.Cells.Item(122, 12).Value = 11048.4999  # L122: was 11098.5
$ws.Cells.Item(122, 13).Value = -6652.3999  # M122: was -6421.75
$ws.Cells.Item(122, 14).Value = -15948.4999  # N122: was -15998.5

$ws.Cells.Item(126, 8).Value = 2429.8096  # H126: was 2643.3333
$ws.Cells.Item(126, 9).Value = 1622.1818  # I126: was 1686.4
$ws.Cells.Item(126, 10).Value = 3318.2  # J126: was 3839.5
$ws.Cells.Item(126, 11).Value = 4866.5454  # K126: was 5059.200000000001
$ws.Cells.Item(126, 12).Value = 9954.599999999999  # L126: was 11518.5
$ws.Cells.Item(126, 13).Value = -2396.5454  # M126: was -2589.200000000001
$ws.Cells.Item(126, 14).Value = -14894.6  # N126: was -16458.5

$ws.Cells.Item(132, 8).Value = 4090  # H132: was 3246
$ws.Cells.Item(132, 9).Value = 3808.5715  # I132: was 2986.2307
$ws.Cells.Item(132, 11).Value = 11425.7145  # K132: was 8958.6921
$ws.Cells.Item(132, 13).Value = -8895.7145  # M132: was -6428.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4441.913  # H132: was 4249.755
$ws.Cells.Item(132, 9).Value = 2330.182  # I132: was 2630.8965
$ws.Cells.Item(132, 10).Value = 9802.462  # J132: was 6597.1
$ws.Cells.Item(132, 11).Value = 6990.545999999999  # K132: was 7892.689499999999
$ws.Cells.Item(132, 12).Value = 29407.386  # L132: was 19791.3
$ws.Cells.Item(132, 13).Value = -4460.545999999999  # M132: was -5362.689499999999
$ws.Cells.Item(132, 14).Value = -34467.386  # N132: was -24851.3

$ws.Cells.Item(138, 8).Value = 50108.43  # H138: was 49628.777
$ws.Cells.Item(138, 10).Value = 50108.43  # J138: was 49628.777
$ws.Cells.Item(138, 12).Value = 50108.43  # L138: was 49628.777
$ws.Cells.Item(138, 14).Value = -60388.43  # N138: was -59908.777

$ws.Cells.Item(139, 8).Value = 65457.5  # H139: was 65555
$ws.Cells.Item(139, 10).Value = 65457.5  # J139: was 65555
$ws.Cells.Item(139, 12).Value = 65457.5  # L139: was 65555
$ws.Cells.Item(139, 14).Value = -75737.5  # N139: was -75835
